$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.726.34"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "3.490.50"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.57"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.99"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.578"
$ws.Range("E8").Value = "  -1.76%  "
$ws.Range("E9").Value = "  +3.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.11"
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "4.093.98"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.16"
$ws.Range("E14").Value = "  +4.54%  "
$ws.Range("D15").Value = "66.719.66"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "3.485.82"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.26"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "389.50"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.89"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.26"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.68"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.06"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  -3.28%  "
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.54"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.92"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("E37").Value = "  -2.28%  "
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.61"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").Value = "2.817.49"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.97"
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0726"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.81"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.43"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "337.64"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.82"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.39"
$ws.Range("E51").Value = "  -0.88%  "
